$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 7667.1665
$ws.Range("J86").Value = 8200.6
$ws.Range("L86").Value = 8200.6
$ws.Range("N86").Value = -10446.6
$ws.Range("H89").Value = 7667.1665
$ws.Range("J89").Value = 8200.6
$ws.Range("L89").Value = 41003
$ws.Range("N89").Value = -52235
$ws.Range("H111").Value = 2717.1428
$ws.Range("I111").Value = 1397.6
$ws.Range("K111").Value = 4192.799999999999
$ws.Range("M111").Value = -1125.799999999999
$ws.Range("H135").Value = 400
$ws.Range("I135").Value = 400
$ws.Range("K135").Value = 3600
$ws.Range("M135").Value = -1065

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2302.625
$ws.Range("I2").Value = 2473
$ws.Range("J2").Value = 1110
$ws.Range("K2").Value = 2473
$ws.Range("L2").Value = 1110
$ws.Range("M2").Value = -2360
$ws.Range("N2").Value = -1336
$ws.Range("H61").Value = 1900
$ws.Range("I61").Value = 1900
$ws.Range("K61").Value = 1900
$ws.Range("M61").Value = -1688
$ws.Range("H63").Value = 10526.25
$ws.Range("I63").Value = 1050
$ws.Range("K63").Value = 1050
$ws.Range("M63").Value = -364
$ws.Range("H66").Value = 10526.25
$ws.Range("I66").Value = 1050
$ws.Range("K66").Value = 5250
$ws.Range("M66").Value = -1818
$ws.Range("H74").Value = 7068
$ws.Range("I74").Value = 6801.5454
$ws.Range("K74").Value = 6801.5454
$ws.Range("M74").Value = -5927.5454
$ws.Range("H77").Value = 7068
$ws.Range("I77").Value = 6801.5454
$ws.Range("K77").Value = 34007.727
$ws.Range("M77").Value = -29639.727
$ws.Range("H116").Value = 2302.625
$ws.Range("I116").Value = 2473
$ws.Range("J116").Value = 1110
$ws.Range("K116").Value = 2473
$ws.Range("L116").Value = 1110
$ws.Range("M116").Value = -179
$ws.Range("N116").Value = -5698
$ws.Range("H122").Value = 3712.0833
$ws.Range("I122").Value = 3354.4
$ws.Range("K122").Value = 10063.2
$ws.Range("M122").Value = -7613.200000000001
$ws.Range("H136").Value = 1900
$ws.Range("I136").Value = 1900
$ws.Range("K136").Value = 5700
$ws.Range("M136").Value = -3150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2302.625
$ws.Range("I3").Value = 2473
$ws.Range("J3").Value = 1110
$ws.Range("K3").Value = 2473
$ws.Range("L3").Value = 1110
$ws.Range("M3").Value = -2359
$ws.Range("N3").Value = -1338
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 925.6667
$ws.Range("I16").Value = 1049.75
$ws.Range("J16").Value = 677.5
$ws.Range("K16").Value = 1049.75
$ws.Range("L16").Value = 677.5
$ws.Range("M16").Value = -762.75
$ws.Range("N16").Value = -1251.5
$ws.Range("H113").Value = 925.6667
$ws.Range("I113").Value = 1049.75
$ws.Range("J113").Value = 677.5
$ws.Range("K113").Value = 1049.75
$ws.Range("L113").Value = 677.5
$ws.Range("M113").Value = 1120.25
$ws.Range("N113").Value = -5017.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 671876.8
$ws.Range("I4").Value = 6506.125
$ws.Range("J4").Value = 3333359.5
$ws.Range("K4").Value = 19518.375
$ws.Range("L4").Value = 10000078.5
$ws.Range("M4").Value = -19406.375
$ws.Range("N4").Value = -10000302.5
$ws.Range("H8").Value = 1334947.1
$ws.Range("I8").Value = 1334947.1
$ws.Range("K8").Value = 4004841.3
$ws.Range("M8").Value = -4004702.3
$ws.Range("H33").Value = 117
$ws.Range("I33").Value = 89.333336
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 536.000016
$ws.Range("L33").Value = 1200
$ws.Range("M33").Value = -253.000016
$ws.Range("N33").Value = -1766
$ws.Range("H64").Value = 750
$ws.Range("I64").Value = 750
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 2250
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -1980
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 750
$ws.Range("I67").Value = 750
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 2250
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -1314
$ws.Range("N67").ClearContents()
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 1319.8
$ws.Range("I132").Value = 1399.75
$ws.Range("K132").Value = 12597.75
$ws.Range("M132").Value = -10067.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9551.25
$ws.Range("I80").Value = 2752.5
$ws.Range("K80").Value = 2752.5
$ws.Range("M80").Value = -1754.5
$ws.Range("H83").Value = 9551.25
$ws.Range("I83").Value = 2752.5
$ws.Range("K83").Value = 13762.5
$ws.Range("M83").Value = -8770.5
$ws.Range("H102").Value = 7429
$ws.Range("I102").Value = 6699.5
$ws.Range("K102").Value = 6699.5
$ws.Range("M102").Value = -5077.5
$ws.Range("H107").Value = 716.1667
$ws.Range("J107").Value = 74.5
$ws.Range("L107").Value = 74.5
$ws.Range("N107").Value = -3914.5
$ws.Range("H126").Value = 6637.3335
$ws.Range("I126").Value = 6637.3335
$ws.Range("K126").Value = 19912.0005
$ws.Range("M126").Value = -17442.0005
$ws.Range("H132").Value = 3777.25
$ws.Range("J132").Value = 4499.3335
$ws.Range("L132").Value = 13498.0005
$ws.Range("N132").Value = -18558.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3324.75
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H68").Value = 3998.6
$ws.Range("I68").Value = 3332.6667
$ws.Range("J68").Value = 4997.5
$ws.Range("K68").Value = 3332.6667
$ws.Range("L68").Value = 4997.5
$ws.Range("M68").Value = -2583.6667
$ws.Range("N68").Value = -6495.5
$ws.Range("H71").Value = 3998.6
$ws.Range("I71").Value = 3332.6667
$ws.Range("J71").Value = 4997.5
$ws.Range("K71").Value = 16663.3335
$ws.Range("L71").Value = 24987.5
$ws.Range("M71").Value = -12919.3335
$ws.Range("N71").Value = -32475.5
$ws.Range("H93").Value = 751.5
$ws.Range("I93").Value = 751.5
$ws.Range("K93").Value = 751.5
$ws.Range("M93").Value = 496.5
$ws.Range("H122").Value = 4964.3335
$ws.Range("I122").Value = 4964.3335
$ws.Range("K122").Value = 14893.0005
$ws.Range("M122").Value = -12443.0005
$ws.Range("H136").Value = 2923.5
$ws.Range("I136").Value = 3055.4285
$ws.Range("K136").Value = 9166.2855
$ws.Range("M136").Value = -6616.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11599.333
$ws.Range("J62").Value = 11599.333
$ws.Range("L62").Value = 11599.333
$ws.Range("N62").Value = -12847.333
$ws.Range("H65").Value = 11599.333
$ws.Range("J65").Value = 11599.333
$ws.Range("L65").Value = 57996.665
$ws.Range("N65").Value = -64236.665
$ws.Range("H100").Value = 2148.5
$ws.Range("I100").Value = 1718.5714
$ws.Range("K100").Value = 3437.1428
$ws.Range("M100").Value = -2896.1428
$ws.Range("H107").Value = 1077.5294
$ws.Range("I107").Value = 904.25
$ws.Range("J107").Value = 1231.5555
$ws.Range("K107").Value = 2712.75
$ws.Range("L107").Value = 3694.6665
$ws.Range("M107").Value = -792.75
$ws.Range("N107").Value = -7534.666499999999
$ws.Range("H113").Value = 1043.3334
$ws.Range("I113").Value = 1138
$ws.Range("K113").Value = 3414
$ws.Range("M113").Value = -1244
$ws.Range("H122").Value = 5114.1816
$ws.Range("I122").Value = 5125.6
$ws.Range("K122").Value = 15376.8
$ws.Range("M122").Value = -12926.8
$ws.Range("H126").Value = 1561.8
$ws.Range("I126").Value = 1483.6
$ws.Range("J126").Value = 1640
$ws.Range("K126").Value = 4450.799999999999
$ws.Range("L126").Value = 4920
$ws.Range("M126").Value = -1980.799999999999
$ws.Range("N126").Value = -9860
